$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.252.23"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "1.819.28"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4669"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3772"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07413"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8713"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.822.30"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.686"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.413"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07080"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008765"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "27.252.46"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.314"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "2.048.85"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.940"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.243"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.310"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08943"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7831"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.180"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.526"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.931"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.096"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01970"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05252"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.271"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5313"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.366"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +20.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.884"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1690"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.597"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5055"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.668"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06333"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.62%  "
